# Resume_shortlisted_candidates.xlsx
# - Fill in the previously-blank "Date" (column C) for existing candidate rows
# - Re-key Total_Experience (H) / Relevent_Experience (I) as text for those rows
# - Insert two new shortlisted candidates (318 - cleveland1, 319 - surmak)
#   before candidate 320, which pushes it from row 12 down to row 14
# - Candidate 320's row also gets its Date filled in and H/I re-keyed as text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# --- 1. Insert two blank rows above the current row 12 (candidate 320) ---
$ws.Range("A12:A13").EntireRow.Insert()

# --- 2. Fill the Date / Experience fields for the pre-existing candidate rows ---
$existingRows = @{
    2  = @{ C = "2024-01-18"; H = "5";  I = "4"  }
    3  = @{ C = "2024-02-20"; H = "10"; I = "9"  }
    4  = @{ C = "2024-02-20"; H = "11"; I = "10" }
    5  = @{ C = "2024-02-20"; H = "12"; I = "11" }
    6  = @{ C = "2024-02-20"; H = "13"; I = "12" }
    7  = @{ C = "2024-02-20"; H = "10"; I = "10" }
    8  = @{ C = "2024-01-02"; H = "16"; I = "15" }
    9  = @{ C = "2024-01-02"; H = "16"; I = "15" }
    10 = @{ C = "2024-01-02"; H = "16"; I = "15" }
    11 = @{ C = "2024-01-02"; H = "16"; I = "15" }
    14 = @{ C = "2024-02-20"; H = "10"; I = "9"  }
}

foreach ($r in $existingRows.Keys) {
    $row = $existingRows[$r]
    Set-TextValue ([int]$r) 3 $row.C
    Set-TextValue ([int]$r) 8 $row.H
    Set-TextValue ([int]$r) 9 $row.I
}

# --- 3. Populate the two newly-inserted candidate rows (318, 319) ---
$newRows = @{
    12 = @{ B = 318; C = "2024-01-02"; D = "java"; E = "cleveland1"; F = "55667788";
            G = "cleveland1@gmail.com"; H = "16"; I = "15"; J = "nasa corporation";
            K = "3"; L = "upgraded for further interview level"; M = "shortlisted" }
    13 = @{ B = 319; C = "2024-01-02"; D = "java"; E = "surmak"; F = "55667788";
            G = "surmka@gmail.com"; H = "16"; I = "15"; J = "umbrala corporation";
            K = "3"; L = "upgraded for further interview level"; M = "shortlisted" }
}

foreach ($r in $newRows.Keys) {
    $row = $newRows[$r]
    $rowNum = [int]$r
    $ws.Cells.Item($rowNum, 2).Value = $row.B
    Set-TextValue $rowNum 3 $row.C
    Set-TextValue $rowNum 4 $row.D
    Set-TextValue $rowNum 5 $row.E
    Set-TextValue $rowNum 6 $row.F
    Set-TextValue $rowNum 7 $row.G
    Set-TextValue $rowNum 8 $row.H
    Set-TextValue $rowNum 9 $row.I
    Set-TextValue $rowNum 10 $row.J
    Set-TextValue $rowNum 11 $row.K
    Set-TextValue $rowNum 12 $row.L
    Set-TextValue $rowNum 13 $row.M
}
